$wb = $excel.ActiveWorkbook

# --- Insert a new "ArchivedSamples" sheet right after "ConclToApprove" ---
# Cloning "DonationInfo" gives us the correct shared style/format/merge
# structure (this new sheet mirrors that layout), then we edit the handful
# of cells that actually differ. Excel places the copy BEFORE the original
# (both keep the same index order), so re-fetch each sheet by name
# afterwards rather than trusting the variable that triggered the copy.
$wb.Worksheets.Item("DonationInfo").Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item("ConclToApprove"))

$newSheet = $wb.Worksheets.Item("DonationInfo (2)")
$newSheet.Name = "ArchivedSamples"
$donationInfo = $wb.Worksheets.Item("DonationInfo")

# Cell content that differs from the DonationInfo template
$newSheet.Range("B3").Value = "ArchiveSamples"
$newSheet.Range("B8").Value = "ArchiveSamples"
$newSheet.Range("B12").Value = "ArchiveSamples"
$newSheet.Range("C12").Value = "55"
$newSheet.Range("A10").Value = "Assert404"

# Row 12 is shorter now ("ArchiveSamples" / "55" vs the longer donation
# strings) so it no longer needs the taller wrapped height.
$newSheet.Rows(12).RowHeight = 30

# Column widths particular to this new sheet
$newSheet.Columns(1).ColumnWidth = 20.71
$newSheet.Columns(2).ColumnWidth = 11.140625
$newSheet.Columns(3).ColumnWidth = 19.86

# Per-sheet selection state (set before the final Activate so it sticks)
$newSheet.Range("I4").Select()
$donationInfo.Range("A1:C13").Select()

# The new sheet is the one left active/selected on open.
$newSheet.Activate()
